# StoryElements.xlsx edit script
# Implements: "added H3, updated H4 & H5, smaller improvements"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / label updates (order matters for shared-string table layout) ---
# R7 used to read "(MainSound(On/Off), FuseSound, metallic Box Sound, torchlight?)";
# shorten it to just "(torchlight?)".
$ws.Range("R7").Value = "(torchlight?)"

# R8 used to read "(silent creepy noise, footsteps, combination with light?)";
# shorten it to just "(combination with light?)".
$ws.Range("R8").Value = "(combination with light?)"

# New note for row 6 (H3 details column)
$ws.Range("R6").Value = "((Message), look at window inpcmode)"

# --- New / updated "%DoneSat" (Q column) percentages ---
$ws.Range("K4").Value = 0.2
$ws.Range("K4").NumberFormat = "0%"

$ws.Range("Q4").Value = 0
$ws.Range("Q4").NumberFormat = "0%"

$ws.Range("Q5").Value = 0
$ws.Range("Q5").NumberFormat = "0%"

$ws.Range("Q6").Value = 0.93
$ws.Range("Q6").NumberFormat = "0%"

$ws.Range("Q7").Value = 0.97
$ws.Range("Q8").Value = 0.97

$ws.Range("Q9").Value = 0
$ws.Range("Q9").NumberFormat = "0%"

$ws.Range("Q10").Value = 0
$ws.Range("Q10").NumberFormat = "0%"

$ws.Range("Q11").Value = 0
$ws.Range("Q11").NumberFormat = "0%"

$ws.Range("Q12").Value = 0
$ws.Range("Q12").NumberFormat = "0%"

$ws.Range("Q13").Value = 0
$ws.Range("Q13").NumberFormat = "0%"

$ws.Range("Q14").Value = 0
$ws.Range("Q14").NumberFormat = "0%"

$ws.Range("Q15").Value = 0
$ws.Range("Q15").NumberFormat = "0%"

$ws.Range("Q16").Value = 0
$ws.Range("Q16").NumberFormat = "0%"

$ws.Range("Q17").Value = 0
$ws.Range("Q17").NumberFormat = "0%"

$ws.Range("Q18").Value = 0
$ws.Range("Q18").NumberFormat = "0%"

$ws.Range("Q19").Value = 0
$ws.Range("Q19").NumberFormat = "0%"

$ws.Range("Q20").Value = 0
$ws.Range("Q20").NumberFormat = "0%"

$ws.Range("Q21").Value = 0
$ws.Range("Q21").NumberFormat = "0%"

$ws.Range("Q22").Value = 0
$ws.Range("Q22").NumberFormat = "0%"

# --- Selection state ---
$ws.Range("N4:N5").Select()
